# Add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计", pushing the
#    existing quarterly sheets one position to the right.
# 2. Populate "2022-Q3" with the new fund-holding detail rows.
# 3. Insert a new row into "总计" (the summary sheet) for "2022-Q3" and
#    shift the previously-existing summary rows down by one.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# --- 1. Create the new sheet right after "总计" ------------------------
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Re-create the look of the other quarter sheets: bold, centered,
# top-aligned cells with a thin border - this single style is shared by
# both the header row (B1:H1) and the index column (A2:A4) on every
# other quarter sheet. (Union ranges aren't reliably supported here, so
# style each area separately.)
foreach ($styledCells in @($q3.Range("B1:H1"), $q3.Range("A2:A4"))) {
    $styledCells.Font.Bold = $true
    $styledCells.HorizontalAlignment = -4108
    $styledCells.VerticalAlignment = -4160
    $styledCells.Borders.LineStyle = 1
}

# --- 2. Populate the new "2022-Q3" sheet --------------------------------
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# B (fund code) and D:G (size/position/value figures) are stored as text
# in the source data (leading zeros, fixed decimal formatting), so force
# a text number-format before writing them - otherwise Excel will infer
# them as numbers and drop things like leading zeros / trailing zeros.
$q3.Range("B2:B4").NumberFormat = "@"
$q3.Range("D2:G4").NumberFormat = "@"

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "013441"
$q3.Cells.Item(2, 3).Value = "西藏东财创新医疗六个月定开混合"
$q3.Cells.Item(2, 4).Value = "0.49"
$q3.Cells.Item(2, 5).Value = "82.53"
$q3.Cells.Item(2, 6).Value = "3.21"
$q3.Cells.Item(2, 7).Value = "0.0157"
$q3.Cells.Item(2, 8).Value = 10

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "005143"
$q3.Cells.Item(3, 3).Value = "中融沪港深大消费主题灵活配置混合C"
$q3.Cells.Item(3, 4).Value = "0.27"
$q3.Cells.Item(3, 5).Value = "90.10"
$q3.Cells.Item(3, 6).Value = "5.75"
$q3.Cells.Item(3, 7).Value = "0.0155"
$q3.Cells.Item(3, 8).Value = 2

$q3.Cells.Item(4, 1).Value = 2
$q3.Cells.Item(4, 2).Value = "005142"
$q3.Cells.Item(4, 3).Value = "中融沪港深大消费主题灵活配置混合A"
$q3.Cells.Item(4, 4).Value = "0.13"
$q3.Cells.Item(4, 5).Value = "90.10"
$q3.Cells.Item(4, 6).Value = "5.75"
$q3.Cells.Item(4, 7).Value = "0.0075"
$q3.Cells.Item(4, 8).Value = 2

# --- 3. Insert a new summary row into "总计" -----------------------------
# Read the existing rows (2..7) first, then shift everything down one row
# (bottom-to-top) before writing in the new 2022-Q3 summary values.
$oldValues = @()
for ($r = 2; $r -le 7; $r++) {
    $oldValues += , @(
        $totalSheet.Cells.Item($r, 1).Value(),
        $totalSheet.Cells.Item($r, 2).Value(),
        $totalSheet.Cells.Item($r, 3).Value(),
        $totalSheet.Cells.Item($r, 4).Value()
    )
}

for ($i = $oldValues.Length - 1; $i -ge 0; $i--) {
    $destRow = $i + 3
    $row = $oldValues[$i]
    $totalSheet.Cells.Item($destRow, 2).Value = $row[1]
    $totalSheet.Cells.Item($destRow, 3).Value = $row[2]
    $totalSheet.Cells.Item($destRow, 4).Value = $row[3]
}

# Row 8 is brand new - copy the bold/centered index-column style from an
# existing styled cell (A7) onto A8 before writing its value.
$totalSheet.Cells.Item(7, 1).Copy()
$totalSheet.Cells.Item(8, 1).PasteSpecial(-4122)
$totalSheet.Cells.Item(8, 1).Value = 6

# Now write the new 2022-Q3 row at row 2 (A2's style/value were already
# correct - index "0" never moved).
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.04
